# Apply NATMI TPM re-run update to Wnt1-Ryk sheet (Sheet1)
# Sending-cluster labels shift (FAPs block -> ECs, MuSCs block -> FAPs) and the
# associated ligand/receptor/edge statistics are refreshed with the new TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("G2").Value = 0.131499
$ws.Range("H2").Value = 0.394497
$ws.Range("I2").Value = 0.3654391092296077
$ws.Range("J2").Value = 0.3654391092296077
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 9.673704333333333
$ws.Range("N2").Value = 29.021113
$ws.Range("O2").Value = 0.1714456165911166
$ws.Range("P2").Value = 0.1714456165911166
$ws.Range("Q2").Value = 1.272082446129
$ws.Range("R2").Value = 11.448742015161
$ws.Range("S2").Value = 0.06265293340837851
$ws.Range("T2").Value = 0.06265293340837851

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("G3").Value = 0.131499
$ws.Range("H3").Value = 0.394497
$ws.Range("I3").Value = 0.3654391092296077
$ws.Range("J3").Value = 0.3654391092296077
$ws.Range("O3").Value = 0.4148961799842911
$ws.Range("P3").Value = 0.4148961799842911
$ws.Range("Q3").Value = 3.078423105927
$ws.Range("R3").Value = 27.705807953343
$ws.Range("S3").Value = 0.1516192904362263
$ws.Range("T3").Value = 0.1516192904362263

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("G4").Value = 0.131499
$ws.Range("H4").Value = 0.394497
$ws.Range("I4").Value = 0.3654391092296077
$ws.Range("J4").Value = 0.3654391092296077
$ws.Range("M4").Value = 23.09142233333333
$ws.Range("N4").Value = 69.27426699999999
$ws.Range("O4").Value = 0.4092458280188166
$ws.Range("P4").Value = 0.4092458280188166
$ws.Range("Q4").Value = 3.036498945411
$ws.Range("R4").Value = 27.328490508699
$ws.Range("S4").Value = 0.1495544308471296
$ws.Range("T4").Value = 0.1495544308471296

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("G5").Value = 0.131499
$ws.Range("H5").Value = 0.394497
$ws.Range("I5").Value = 0.3654391092296077
$ws.Range("J5").Value = 0.3654391092296077
$ws.Range("M5").Value = 0.2489653333333333
$ws.Range("N5").Value = 0.746896
$ws.Range("O5").Value = 0.004412375405775742
$ws.Range("P5").Value = 0.004412375405775742
$ws.Range("Q5").Value = 0.03273869236800001
$ws.Range("R5").Value = 0.294648231312
$ws.Range("S5").Value = 0.001612454537873316
$ws.Range("T5").Value = 0.001612454537873316

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("G6").Value = 0.2283393333333333
$ws.Range("H6").Value = 0.685018
$ws.Range("I6").Value = 0.6345608907703922
$ws.Range("J6").Value = 0.6345608907703922
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 9.673704333333333
$ws.Range("N6").Value = 29.021113
$ws.Range("O6").Value = 0.1714456165911166
$ws.Range("P6").Value = 0.1714456165911166
$ws.Range("Q6").Value = 2.208887198337111
$ws.Range("R6").Value = 19.879984785034
$ws.Range("S6").Value = 0.1087926831827381
$ws.Range("T6").Value = 0.1087926831827381

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("G7").Value = 0.2283393333333333
$ws.Range("H7").Value = 0.685018
$ws.Range("I7").Value = 0.6345608907703922
$ws.Range("J7").Value = 0.6345608907703922
$ws.Range("O7").Value = 0.4148961799842911
$ws.Range("P7").Value = 0.4148961799842911
$ws.Range("Q7").Value = 5.345478518660222
$ws.Range("R7").Value = 48.10930666794199
$ws.Range("S7").Value = 0.2632768895480647
$ws.Range("T7").Value = 0.2632768895480647

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("G8").Value = 0.2283393333333333
$ws.Range("H8").Value = 0.685018
$ws.Range("I8").Value = 0.6345608907703922
$ws.Range("J8").Value = 0.6345608907703922
$ws.Range("M8").Value = 23.09142233333333
$ws.Range("N8").Value = 69.27426699999999
$ws.Range("O8").Value = 0.4092458280188166
$ws.Range("P8").Value = 0.4092458280188166
$ws.Range("Q8").Value = 5.272679981311777
$ws.Range("R8").Value = 47.454119831806
$ws.Range("S8").Value = 0.259691397171687
$ws.Range("T8").Value = 0.259691397171687

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("G9").Value = 0.2283393333333333
$ws.Range("H9").Value = 0.685018
$ws.Range("I9").Value = 0.6345608907703922
$ws.Range("J9").Value = 0.6345608907703922
$ws.Range("M9").Value = 0.2489653333333333
$ws.Range("N9").Value = 0.746896
$ws.Range("O9").Value = 0.004412375405775742
$ws.Range("P9").Value = 0.004412375405775742
$ws.Range("Q9").Value = 0.05684857823644445
$ws.Range("R9").Value = 0.5116372041280001
$ws.Range("S9").Value = 0.002799920867902426
$ws.Range("T9").Value = 0.002799920867902426

